# Update SCD0026-001 until SCD0026-017 Fix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from SCD0338 to SCD0026
$ws.Name = "SCD0026"

# Update the TC_ID values in column B (rows 2-6) from SCD0338-016 to SCD0026-016
$ws.Range("B2:B6").Value = "SCD0026-016"

# Move the active selection from G5 to B7
$ws.Range("B7").Select()
